# heat_exchanger_input1.xlsx : "Excel parameter reading/input is moved to
# heat_exchanger.py, excel-input file modified."
#
# Adds new "hex_name" / "no_points" / "calc_type" parameters to the
# Geometry sheet and a new "units" parameter to both Fluid_1 and Fluid_2
# sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "Geometry": three new rows (14-16) ---
$wsGeom = $wb.Worksheets.Item("Geometry")
$wsGeom.Range("A14").Value = "hex_name"
$wsGeom.Range("B14").Value = "condenser0"
$wsGeom.Range("A15").Value = "no_points"
$wsGeom.Range("B15").Value = 100
$wsGeom.Range("A16").Value = "calc_type"
$wsGeom.Range("B16").Value = "const"

# --- Sheet "Fluid_1": one new row (9) with the "units" parameter ---
$wsFluid1 = $wb.Worksheets.Item("Fluid_1")
$wsFluid1.Range("A9").Value = "units"
$wsFluid1.Range("B9").Value = 21

# --- Sheet "Fluid_2": fill in the previously blank row 8 with "units" ---
$wsFluid2 = $wb.Worksheets.Item("Fluid_2")
$wsFluid2.Range("A8").Value = "units"
$wsFluid2.Range("B8").Value = 21

# Column A got a little wider to fit the longer labels (closest width the
# engine's column-width quantization can reach to the authored 19.6328125)
$wsFluid2.Columns("A").ColumnWidth = 18.8

# --- Sheet "Problem_description": only the selected cell moved ---
$wsProblem = $wb.Worksheets.Item("Problem_description")

# Update the on-screen selection for every sheet to match where the user
# last clicked while editing. "Fluid_2" is the sheet that stays active /
# selected in the saved workbook, so select it last.
$wsGeom.Activate()
$wsGeom.Range("A17").Select()

$wsFluid1.Activate()
$wsFluid1.Range("B9").Select()

$wsProblem.Activate()
$wsProblem.Range("E12").Select()

$wsFluid2.Activate()
$wsFluid2.Range("A9").Select()
